$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Top Scores" sheet. Worksheets.Add() (with no arguments)
#    inserts the new sheet before the currently active sheet, i.e. at the
#    very front of the tab strip - exactly where it needs to land.
# ---------------------------------------------------------------------------
$top = $wb.Worksheets.Add()
$top.Name = "Top Scores"

$classic = $wb.Worksheets.Item("Classic Scores")
$hangman = $wb.Worksheets.Item("Hangman Scores")

# ---------------------------------------------------------------------------
# 2. "Top Scores" - combined Game mode / Difficulty / Score / Time table.
# ---------------------------------------------------------------------------
$top.Cells.Item(1,1).Value = "Game mode"
$top.Cells.Item(1,2).Value = "Difficulty"
$top.Cells.Item(1,3).Value = "Score"
$top.Cells.Item(1,4).Value = "Time"

$top.Cells.Item(2,1).Value = "Classic"
$top.Cells.Item(2,2).Value = """Normal"""
$top.Cells.Item(2,3).Value = 40
$top.Cells.Item(2,4).Value = 20

$top.Cells.Item(3,1).Value = "Classic"
$top.Cells.Item(3,2).Value = """Easy"""
$top.Cells.Item(3,3).Value = 40
$top.Cells.Item(3,4).Value = 20

$top.Cells.Item(4,1).Value = "Classic"
$top.Cells.Item(4,2).Value = """Hard"""
$top.Cells.Item(4,3).Value = 20
$top.Cells.Item(4,4).Value = 10

$top.Cells.Item(5,1).Value = "Classic"
$top.Cells.Item(5,2).Value = """Ultra_Hard"""
$top.Cells.Item(5,3).Value = 80
$top.Cells.Item(5,4).Value = 50

$top.Cells.Item(6,1).Value = "Hangman"
$top.Cells.Item(6,2).Value = """Normal"""
$top.Cells.Item(6,3).Value = 40
$top.Cells.Item(6,4).Value = 20

$top.Cells.Item(7,1).Value = "Hangman"
$top.Cells.Item(7,2).Value = """Easy"""
$top.Cells.Item(7,3).Value = 40
$top.Cells.Item(7,4).Value = 20

$top.Cells.Item(8,1).Value = "Hangman"
$top.Cells.Item(8,2).Value = """Hard"""
$top.Cells.Item(8,3).Value = 20
$top.Cells.Item(8,4).Value = 10

$top.Cells.Item(9,1).Value = "Hangman"
$top.Cells.Item(9,2).Value = """Ultra_Hard"""
$top.Cells.Item(9,3).Value = 80
$top.Cells.Item(9,4).Value = 50

$top.Columns.Item(1).ColumnWidth = 10.28

# ---------------------------------------------------------------------------
# 3. "Classic Scores" - Difficulty / Score / Time table, becomes the sheet
#    that is active when the workbook is reopened.
# ---------------------------------------------------------------------------
$classic.Cells.Item(1,1).Value = "Difficulty"
$classic.Cells.Item(1,2).Value = "Score"
$classic.Cells.Item(1,3).Value = "Time"

$classic.Cells.Item(2,1).Value = """Normal"""
$classic.Cells.Item(2,2).Value = 40
$classic.Cells.Item(2,3).Value = 20

$classic.Cells.Item(3,1).Value = """Normal"""
$classic.Cells.Item(3,2).Value = 20
$classic.Cells.Item(3,3).Value = 10

$classic.Cells.Item(4,1).Value = """Easy"""
$classic.Cells.Item(4,2).Value = 40
$classic.Cells.Item(4,3).Value = 20

$classic.Cells.Item(5,1).Value = """Hard"""
$classic.Cells.Item(5,2).Value = 20
$classic.Cells.Item(5,3).Value = 10

$classic.Cells.Item(6,1).Value = """Ultra_Hard"""
$classic.Cells.Item(6,2).Value = 80
$classic.Cells.Item(6,3).Value = 50

# ---------------------------------------------------------------------------
# 4. "Hangman Scores" - same Difficulty / Score / Time layout.
# ---------------------------------------------------------------------------
$hangman.Cells.Item(1,1).Value = "Difficulty"
$hangman.Cells.Item(1,2).Value = "Score"
$hangman.Cells.Item(1,3).Value = "Time"

$hangman.Cells.Item(2,1).Value = """Normal"""
$hangman.Cells.Item(2,2).Value = 40
$hangman.Cells.Item(2,3).Value = 20

$hangman.Cells.Item(3,1).Value = """Normal"""
$hangman.Cells.Item(3,2).Value = 20
$hangman.Cells.Item(3,3).Value = 10

$hangman.Cells.Item(4,1).Value = """Easy"""
$hangman.Cells.Item(4,2).Value = 40
$hangman.Cells.Item(4,3).Value = 20

$hangman.Cells.Item(5,1).Value = """Hard"""
$hangman.Cells.Item(5,2).Value = 20
$hangman.Cells.Item(5,3).Value = 10

$hangman.Cells.Item(6,1).Value = """Ultra_Hard"""
$hangman.Cells.Item(6,2).Value = 80
$hangman.Cells.Item(6,3).Value = 50

# ---------------------------------------------------------------------------
# 5. Selections per-sheet, then activate "Classic Scores" last so it ends up
#    as the active tab (activeTab) when the file is saved.
# ---------------------------------------------------------------------------
$hangman.Range("A1:C6").Select() | Out-Null
$top.Range("G10").Select() | Out-Null
$classic.Range("M10").Select() | Out-Null
$classic.Activate()
